# Adds interrupt-related columns (GPIO Pull-up/pull-down, NVIC interrupt table)
# to the pin description sheet for pins PC1, PA8, PA9, bolds the header row,
# and tweaks page setup / selection to match the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column E: "GPIO Pull-up/pull-down" -------------------------------
# Fill column E completely (header + data) before moving on to column F so
# new shared-string entries are interleaved in the same order as the
# original authoring session (E-header, E-data, then F-header, F-data).
$ws.Range("E1").Value = "GPIO Pull-up/pull-down"
$ws.Range("E10").Value = "pull-up"
$ws.Range("E11").Value = "pull-up"

# --- New column F: "NVIC interrupt table" ----------------------------------
$ws.Range("F1").Value = "NVIC interrupt table"
$ws.Range("F10").Value = "enabled"
$ws.Range("F11").Value = "enabled"
$ws.Range("F12").Value = "enabled"

# --- Header row formatting --------------------------------------------------
$ws.Range("A1:F1").Font.Bold = $true

# --- Column widths for the two new columns ---------------------------------
$ws.Columns.Item(5).ColumnWidth = 22.333333333333332
$ws.Columns.Item(6).ColumnWidth = 23

# --- Page setup --------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection as left by the author ----------------------------------------
$ws.Range("B12").Select()
